$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1108
$ws.Range("J19").Value = 1080.8334
$ws.Range("L19").Value = 1080.8334
$ws.Range("N19").Value = -1430.8334

$ws.Range("H47").Value = 5533.5
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()

$ws.Range("H87").Value = 21463.223
$ws.Range("J87").Value = 21463.223
$ws.Range("L87").Value = 21463.223
$ws.Range("N87").Value = -23959.223

$ws.Range("H90").Value = 21463.223
$ws.Range("J90").Value = 21463.223
$ws.Range("L90").Value = 64389.66900000001
$ws.Range("N90").Value = -76869.66900000001

$ws.Range("H103").Value = 2111.5557
$ws.Range("I103").Value = 1001
$ws.Range("K103").Value = 3003
$ws.Range("M103").Value = -2417

$ws.Range("H114").Value = 38666
$ws.Range("J114").Value = 38666
$ws.Range("L114").Value = 38666
$ws.Range("N114").Value = -47344

$ws.Range("H132").Value = 25370.334
$ws.Range("I132").Value = 1339.8235
$ws.Range("K132").Value = 4019.4705
$ws.Range("M132").Value = -1489.4705

$ws.Range("H137").Value = 2264.9355
$ws.Range("I137").Value = 1695.8
$ws.Range("J137").Value = 2798.5
$ws.Range("K137").Value = 5087.4
$ws.Range("L137").Value = 8395.5
$ws.Range("M137").Value = -2537.4
$ws.Range("N137").Value = -13495.5

$ws.Range("H138").Value = 3666.6667
$ws.Range("J138").Value = 6000
$ws.Range("L138").Value = 18000
$ws.Range("N138").Value = -28280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3532.44
$ws.Range("J61").Value = 6930
$ws.Range("L61").Value = 6930
$ws.Range("N61").Value = -7354

$ws.Range("H88").Value = 1575.7273
$ws.Range("I88").Value = 1120
$ws.Range("J88").Value = 1836.1428
$ws.Range("K88").Value = 1120
$ws.Range("L88").Value = 1836.1428
$ws.Range("M88").Value = -714
$ws.Range("N88").Value = -2648.1428

$ws.Range("H91").Value = 1575.7273
$ws.Range("I91").Value = 1120
$ws.Range("J91").Value = 1836.1428
$ws.Range("K91").Value = 1120
$ws.Range("L91").Value = 1836.1428
$ws.Range("M91").Value = 284
$ws.Range("N91").Value = -4644.1428

$ws.Range("H132").Value = 3597.913
$ws.Range("I132").Value = 2546.125
$ws.Range("J132").Value = 6002
$ws.Range("K132").Value = 7638.375
$ws.Range("L132").Value = 18006
$ws.Range("M132").Value = -5108.375
$ws.Range("N132").Value = -23066

$ws.Range("H136").Value = 3532.44
$ws.Range("J136").Value = 6930
$ws.Range("L136").Value = 20790
$ws.Range("N136").Value = -25890

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H130").Value = 50042.89
$ws.Range("J130").Value = 50042.89
$ws.Range("L130").Value = 50042.89
$ws.Range("N130").Value = -60082.89

$ws.Range("H134").Value = 2463.2334
$ws.Range("I134").Value = 1664.625
$ws.Range("K134").Value = 4993.875
$ws.Range("M134").Value = -2458.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3341.795
$ws.Range("I31").Value = 2231.04
$ws.Range("J31").Value = 5325.2856
$ws.Range("K31").Value = 2231.04
$ws.Range("L31").Value = 5325.2856
$ws.Range("M31").Value = -1936.04
$ws.Range("N31").Value = -5915.2856

$ws.Range("H33").Value = 1343.6666
$ws.Range("I33").Value = 1343.6666
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 1343.6666
$ws.Range("L33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("N33").Value = -964.6666

$ws.Range("H34").Value = 3341.795
$ws.Range("I34").Value = 2231.04
$ws.Range("J34").Value = 5325.2856
$ws.Range("K34").Value = 2231.04
$ws.Range("L34").Value = 5325.2856
$ws.Range("M34").Value = -2029.04
$ws.Range("N34").Value = -5729.2856

$ws.Range("H58").Value = 4335.636
$ws.Range("I58").Value = 2928.6
$ws.Range("K58").Value = 2928.6
$ws.Range("M58").Value = -2725.6

$ws.Range("H99").Value = 14967751
$ws.Range("I99").Value = 4884649
$ws.Range("K99").Value = 4884649
$ws.Range("M99").Value = -4883151

$ws.Range("H122").Value = 368087.25
$ws.Range("I122").Value = 786160.25
$ws.Range("K122").Value = 2358480.75
$ws.Range("M122").Value = -2356030.75

$ws.Range("H126").Value = 14967751
$ws.Range("I126").Value = 4884649
$ws.Range("K126").Value = 14653947
$ws.Range("M126").Value = -14651477

$ws.Range("H132").Value = 4112.25
$ws.Range("I132").Value = 4112.25
$ws.Range("K132").Value = 12336.75
$ws.Range("M132").Value = -9806.75

$ws.Range("H134").Value = 4909.0713
$ws.Range("I134").Value = 3528.0417
$ws.Range("K134").Value = 10584.1251
$ws.Range("M134").Value = -8049.125100000001

$ws.Range("H136").Value = 4335.636
$ws.Range("I136").Value = 2928.6
$ws.Range("K136").Value = 8785.799999999999
$ws.Range("M136").Value = -6235.799999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 122.166664
$ws.Range("I2").Value = 21
$ws.Range("K2").Value = 126
$ws.Range("M2").Value = -13

$ws.Range("H17").Value = 1033.3334
$ws.Range("I17").Value = 1000
$ws.Range("K17").Value = 3000
$ws.Range("M17").Value = -2831

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2944.0667
$ws.Range("I102").Value = 1795.375
$ws.Range("J102").Value = 4256.857
$ws.Range("K102").Value = 1795.375
$ws.Range("L102").Value = 4256.857
$ws.Range("M102").Value = -173.375
$ws.Range("N102").Value = -7500.857

$ws.Range("H131").Value = 103672.57
$ws.Range("J131").Value = 119601.78
$ws.Range("L131").Value = 119601.78
$ws.Range("N131").Value = -129681.78

$ws.Range("H135").Value = 60526.26
$ws.Range("J135").Value = 60526.26
$ws.Range("L135").Value = 60526.26
$ws.Range("N135").Value = -70666.26000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5224.0527
$ws.Range("I132").Value = 3590.889
$ws.Range("J132").Value = 9232.727999999999
$ws.Range("K132").Value = 10772.667
$ws.Range("L132").Value = 27698.184
$ws.Range("M132").Value = -8242.667000000001
$ws.Range("N132").Value = -32758.184

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2170.5881
$ws.Range("I81").Value = 2109.2144
$ws.Range("J81").Value = 2457
$ws.Range("K81").Value = 4218.4288
$ws.Range("L81").Value = 4914
$ws.Range("M81").Value = -3157.4288
$ws.Range("N81").Value = -7036

$ws.Range("H84").Value = 2170.5881
$ws.Range("I84").Value = 2109.2144
$ws.Range("J84").Value = 2457
$ws.Range("K84").Value = 21092.144
$ws.Range("L84").Value = 24570
$ws.Range("M84").Value = -15788.144
$ws.Range("N84").Value = -35178

$ws.Range("H131").Value = 58929.168
$ws.Range("J131").Value = 58929.168
$ws.Range("L131").Value = 58929.168
$ws.Range("N131").Value = -69009.16800000001

$ws.Range("H132").Value = 1667.4348
$ws.Range("I132").Value = 1024.6333
$ws.Range("J132").Value = 2872.6875
$ws.Range("K132").Value = 3073.8999
$ws.Range("L132").Value = 8618.0625
$ws.Range("M132").Value = -543.8998999999999
$ws.Range("N132").Value = -13678.0625

$ws.Range("H136").Value = 2048.842
$ws.Range("I136").Value = 1143.5667
$ws.Range("K136").Value = 3430.7001
$ws.Range("M136").Value = -880.7001
